# The diff for this commit touches only two kinds of things, neither of
# which is reachable (or meaningful) through the Word UI / COM object
# model:
#
#   1. The order of the xmlns:* attributes on the root element of
#      word/document.xml, word/footer.xml, word/header.xml and
#      word/styles.xml is shuffled. This is purely a side effect of the
#      XML serializer used by the tool that produced the new package
#      (a different JAXB/JVM run, per the comment below) - it carries no
#      document content, formatting, or structure, and Word's COM
#      automation surface has no property that controls namespace
#      *declaration order* on save.
#
#   2. A single word changes inside an XML processing comment that
#      docx4j stamps as the very first child of <w:body>:
#
#        <!-- Modified by docx4j 11.5.6 (Apache licensed) using
#             REFERENCE JAXB in Oracle Java 21.0.8 on Mac OS X -->
#                                  ^^^^^^
#      becomes "... in Microsoft Java 21.0.8 ...". That is a raw XML
#      comment left behind by the docx4j build tool itself (recording
#      which JVM produced the file) - it is not Word content, not a
#      Word "comment" (no w:comment/w:commentReference involved, no
#      comments.xml), and it is not part of the WordprocessingML object
#      model at all, so there is no Range/Selection/Document API that
#      can read or write it.
#
# In other words, nothing a user could do inside Word (type text, change
# formatting, touch a property exposed on Application/Document/Range/...)
# corresponds to this diff; it is tooling/environment provenance noise,
# not a document edit. The safest - and closest - match is therefore to
# leave the document completely untouched: any real mutation would cause
# the host to re-serialize the edited parts from the object model, which
# drops this non-model comment entirely (making the result further from
# the target, not closer) while still being unable to reproduce the
# vendor-string text or the namespace ordering.
#
# So: touch nothing.
$d = $word.ActiveDocument
